$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = -0
    "B2" = -0.107670151431414
    "C2" = -0
    "D2" = 0.2053876115891933
    "E2" = -0.000870470479885202
    "G2" = 0
    "J2" = -0
    "K2" = -0.01825145604069086
    "L2" = -0
    "M2" = 0.2064948814758124
    "N2" = 0.00208025141448778
    "R2" = -0
    "S2" = 0
    "T2" = -0.09778338212392844
    "V2" = 0.01432430740601758
    "W2" = -0.04174952113749071
    "AB2" = 0
    "AC2" = -0.07844340024108128
    "AD2" = 0
    "AE2" = -0.01123724989979894
    "AF2" = 0.00467995902030148
    "AG2" = -0
    "AI2" = -0
    "AJ2" = 0
    "AK2" = -0
    "AL2" = -0.03477317437924476
    "AM2" = 0
    "AN2" = 0.02975591359928845
    "AO2" = 0.06894031625633901
    "AR2" = -0
    "AT2" = 0
    "AU2" = -0.1504158949591624
    "AW2" = 0.06762809990804736
    "AX2" = 0.01316411438484101
    "AY2" = -0
    "BC2" = -0
    "BD2" = -0.03736677967854636
    "BF2" = 0.07026386262018187
    "BG2" = 0.02677008819260617
    "BI2" = -0
    "BJ2" = -0
    "BL2" = 0
    "BM2" = 0.01906614199068762
    "BO2" = -0.05327664047297926
    "BP2" = -0.06915173979540963
    "BU2" = 0
    "BV2" = -0.01804737642947697
    "BX2" = 0.01512898727912453
    "BY2" = -0.005352203034577101
    "BZ2" = -0
    "CB2" = 0
    "CD2" = 0
    "CE2" = 0.0336355770062513
    "CG2" = -0.03569317361011079
    "CH2" = 0.0251507005628732
    "CJ2" = -0
    "CM2" = -0
    "CN2" = -0.01586600892987506
    "CO2" = -0
    "CP2" = 0.02594067695371425
    "CQ2" = 0.05198710927025178
    "CT2" = 0
    "CU2" = -0
    "CV2" = -0
    "CW2" = 0.04516854268119422
    "CY2" = -0.03815737946181457
    "CZ2" = 0.01954564062428951
    "DB2" = 0
    "DE2" = -0
    "DF2" = 0.04358123142453846
    "DH2" = 0.01717012150273894
    "DI2" = 0.04988333764285281
    "DJ2" = 0
    "DK2" = -0
    "DL2" = -0
    "DN2" = 0
    "DO2" = -0.02121619411584265
    "DQ2" = 0.03359901195050785
    "DR2" = -0.01385797544263802
    "DS2" = -0
    "DW2" = 0
    "DX2" = -0.05212903656296591
    "DY2" = -0
    "DZ2" = -0.004971154324998677
    "EA2" = -0.02718932500684638
    "EB2" = 0
    "EF2" = -0
    "EG2" = 0.05015491900709527
    "EI2" = 0.07613886383889099
    "EJ2" = -0.0199689068547195
    "EO2" = 0
    "EP2" = 0.04747992163333337
    "EQ2" = 0
    "ER2" = -0.03421357883596667
    "ES2" = 0.04056784176355664
    "ET2" = 0
    "EU2" = -0
    "EV2" = 0
    "EX2" = 0
    "EY2" = 0.047632734695192
    "FA2" = -0.03127982317017986
    "FB2" = 0.02441504420193668
    "FD2" = -0
    "FF2" = -0
    "FG2" = -0
    "FH2" = -0.0002063141115742818
    "FI2" = 0
    "FJ2" = -0.01377775761390556
    "FK2" = -0.001886070509032845
    "FL2" = -0
    "FP2" = -0
    "FQ2" = -0.01316388935840617
    "FR2" = -0
    "FS2" = -0.0009316850821103231
    "FT2" = 0.0192548302470746
    "FV2" = -0
    "FW2" = -0
    "FY2" = 0
    "FZ2" = -0.03444598377532641
    "GA2" = -0
    "GB2" = 0.04121121055052257
    "GD2" = 0
    "GE2" = -0
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

Write-Output "Applied $($values.Count) cell updates to row 2"